# This script applies the "L1BridgeRegistry" gas-log refactor described in the
# commit: renames L2Registry -> L1BridgeRegistry, inserts two new transaction
# rows (DepositManager.approveAndCall "deposit TON to DAO Candidate" and a
# second Layer2Contract.updateSeigniorage "with operator's staking" entry),
# and refreshes the recorded GasUsed figures for many rows.
#
# Rather than trying to replicate the row-shuffling/insert sequence performed
# interactively in Excel, we simply rewrite the whole data table (A1:F60) to
# its final, known-correct contents in one shot -- this is equivalent in
# effect and far less error-prone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final contents of every row of the table, column order A,B,C,D,E,F.
$data = @(
    @('Contract', 'Name', 'Tx', 'Description', 'GasUsed', 'GasUsed (ETH) '),
    @('L1BridgeRegistry', 'addManager', '', '', '50991', '0.00050991'),
    @('OperatorFactory', 'setAddresses', '', '', '115592', '0.00115592'),
    @('Layer2Manager', 'setAddresses', '', '', '217793', '0.00217793'),
    @('Layer2Manager', 'setMinimumInitialDepositAmount', '', '', '57293', '0.00057293'),
    @('SystemConfig', 'setAddresses', '', '', '146738', '0.00146738'),
    @('L1BridgeRegistry', 'registerSystemConfigByManager', '', '', '91010', '0.0009101'),
    @('L1BridgeRegistry', 'setAddresses', '', '', '146750', '0.0014675'),
    @('L1BridgeRegistry', 'registerSystemConfigByManager', '', 'Thanos SystemConfig ', '100363', '0.00100363'),
    @('Layer2Manager', 'registerLayer2Candidate', '', 'using registerLayer2Candidate function', '4862386', '0.04862386'),
    @('SeigManagerV1_3', 'setLayer2StartBlock', '', '', '56913', '0.00056913'),
    @('DepositManager', 'approveAndCall', '', 'deposit TON to Titan Candidate', '374814', '0.00374814'),
    @('DepositManager', 'deposit(address,uint256)', '', 'deposit WTON to Titan Candidate', '303991', '0.00303991'),
    @('DepositManager', 'deposit(address,address,uint256)', '', 'deposit WTON to Titan Candidate', '254169', '0.00254169'),
    @('SeigManagerV1_3', 'updateSeigniorageLayer(layerAddress)', '', 'the first updateSeigniorage : no give seigniorage to l2', '429683', '0.00429683'),
    @('SeigManagerV1_3', 'updateSeigniorageLayer(layerAddress)', '', 'the second updateSeigniorage : not operator', '437530', '0.0043753'),
    @('SeigManagerV1_3', 'updateSeigniorage(1)', '', 'the third updateSeigniorage of operator with claiming', '482606', '0.00482606'),
    @('SeigManagerV1_3', 'updateSeigniorage(2)', '', 'the forth updateSeigniorage of operator with staking ', '588145', '0.00588145'),
    @('DepositManager', 'requestWithdrawal', '', '', '340683', '0.00340683'),
    @('DepositManager', 'processRequest', '', '', '197362', '0.00197362'),
    @('DepositManager', 'approveAndCall', '', 'deposit TON to DAOCandidate', '315110', '0.0031511'),
    @('DepositManager', 'approveAndCall', '', 'deposit TON to DAOCandidate', '356170', '0.0035617'),
    @('DepositManager', 'deposit(address,uint256)', '', 'deposit WTON to DAOCandidate', '249099', '0.00249099'),
    @('DepositManager', 'deposit(address,address,uint256)', '', 'deposit WTON to DAOCandidate', '290666', '0.00290666'),
    @('SeigManagerV1_3', 'updateSeigniorageLayer(layerAddress)', '', 'updateSeigniorage to DAOCandidate', '347933', '0.00347933'),
    @('DepositManager', 'requestWithdrawal(address,uint256)', '', '', '278176', '0.00278176'),
    @('DepositManager', 'processRequests(address,uint256,bool)', '', '', '145889', '0.00145889'),
    @('DepositManager', 'approveAndCall', '', 'deposit TON to Titan Layer2 Candidate', '357714', '0.00357714'),
    @('DepositManager', 'approveAndCall', '', 'deposit TON to DAO Candidate', '322941', '0.00322941'),
    @('DepositManager', 'withdrawAndDepositL2', '', '', '480828', '0.00480828'),
    @('L1BridgeRegistry', 'rejectLayer2Candidate', '', 'legacySystemConfig', '115789', '0.00115789'),
    @('DepositManager', 'approveAndCall', '', 'deposit TON to DAO Candidate', '340614', '0.00340614'),
    @('L1BridgeRegistry', 'registerLayer2Candidate', '', 'thanos SystemConfig', '4823495', '0.04823495'),
    @('SeigManager', 'updateSeigniorageLayer', '', 'titanLayerAddress', '355945', '0.00355945'),
    @('titanLayerContract', 'updateSeigniorage', '', 'with claim', '355880', '0.0035588'),
    @('titanLayerContract', 'updateSeigniorage', '', 'with staking', '355880', '0.0035588'),
    @('DepositManager', 'deposit(address,uint256)', '', '', '286891', '0.00286891'),
    @('SeigManager', 'updateSeigniorageLayer', '', 'no give seigniorage to l2', '435149', '0.00435149'),
    @('DepositManager', 'deposit(address,address,uint256)', '', '', '254169', '0.00254169'),
    @('SeigManager', 'updateSeigniorageLayer', '', 'give seigniorage to l2', '408534', '0.00408534'),
    @('SeigManager', 'updateSeigniorageLayer', '', 'not operator', '408534', '0.00408534'),
    @('Layer2Contract', 'updateSeigniorage', '', 'operator', '487810', '0.0048781'),
    @('Layer2Contract', 'updateSeigniorage', '', 'with operator''s staking', '593349', '0.00593349'),
    @('DepositManager', 'requestWithdrawal', '', '', '341275', '0.00341275'),
    @('DepositManager', 'processRequest', '', '', '180262', '0.00180262'),
    @('L1BridgeRegistry', 'restoreLayer2Candidate', '', '', '72294', '0.00072294'),
    @('SeigManager', 'updateSeigniorageLayer', '', '', '426030', '0.0042603'),
    @('SeigManager', 'updateSeigniorageLayer', '', '', '403330', '0.0040333'),
    @('SeigManager', 'updateSeigniorageLayer', '', '', '403330', '0.0040333'),
    @('Layer2Contract', 'updateSeigniorage', '', 'with operator''s claim', '465506', '0.00465506'),
    @('Layer2Contract', 'updateSeigniorage', '', 'with operator''s staking', '567828', '0.00567828'),
    @('Layer2Contract', 'updateSeigniorage', '', 'with operator''s staking', '437710', '0.0043771'),
    @('DepositManager', 'requestWithdrawal', '', '', '324175', '0.00324175'),
    @('DepositManager', 'processRequest', '', '', '128962', '0.00128962'),
    @('TonContract', 'approveAndCall', '', 'DepositManager.onApprove', '322929', '0.00322929'),
    @('DepositManager', 'deposit(address,uint256)', '', '', '249099', '0.00249099'),
    @('DepositManager', 'deposit(address,address,uint256)', '', '', '257425', '0.00257425'),
    @('SeigManager', 'updateSeigniorageLayer', '', '', '354699', '0.00354699'),
    @('DepositManager', 'requestWithdrawal', '', '', '326832', '0.00326832'),
    @('DepositManager', 'processRequests', '', '', '109935', '0.00109935')
)

$rowCount = $data.Count
$colCount = 6

# Remove any pre-existing content/formatting below the table so no stray
# rows linger past row 60, then make sure the sheet is otherwise empty.
$usedRange = $ws.UsedRange
$ws.Cells.Clear() | Out-Null

# Build a proper 2-D array (PowerShell's array-of-arrays does not marshal
# correctly to a SAFEARRAY when assigned straight to Range.Value).
$arr = New-Object 'object[,]' $rowCount, $colCount
for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $arr[$r, $c] = $data[$r][$c]
    }
}

$targetRange = $ws.Range("A1").Resize($rowCount, $colCount)

# All values in this sheet are stored as text (even the numeric-looking gas
# amounts), matching the source workbook's t="str" cells, so force the
# number format to Text before writing the values.
$targetRange.NumberFormat = "@"
$targetRange.Value = $arr

$ws.Range("A1:F" + $rowCount).EntireColumn.AutoFit() | Out-Null

Write-Host "Updated table to $rowCount rows"
